{"js": "// The document contains a paragraph holding a Word FIELD (fldChar begin /\n// instrText tokens for \" m:'doc.html'.fromHTMLURI() \" / fldChar end).\n// The edit rewrites that field into plain literal text runs that spell out\n// the same token stream wrapped in literal \"{\" and \"}\" characters instead of\n// the field delimiters (i.e. \"{m:'doc.html'.fromHTMLURI()}\"), while leaving\n// the bookmark (\"_GoBack\") that sits between the \"doc.html\" and\n// \"'.fromHTMLURI()\" tokens untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Find the paragraph that contains the field (fldChar begin/end +\n// instrText runs). We detect it through its OOXML rather than assuming a\n// fixed index, so the script is resilient to minor document changes.\nlet targetParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const ooxmlResult = paragraphs.items[i].getOoxml();\n  paragraphs.items[i].__ooxmlResult = ooxmlResult;\n}\nawait context.sync();\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const xml = paragraphs.items[i].__ooxmlResult.value;\n  if (xml && xml.indexOf(\"instrText\") !== -1 && xml.indexOf(\"fldChar\") !== -1) {\n    targetParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the paragraph containing the field to rewrite.\");\n}\n\n// Recover the original <w:p ...> opening tag attributes (rsid bookkeeping\n// etc.) so the rewritten paragraph keeps them instead of becoming a bare\n// <w:p>. The w14:paraId/w14:textId attributes are synthesized by the OOXML\n// exporter itself (not present in the source document), so they are\n// stripped back out.\nconst targetXml = targetParagraph.__ooxmlResult.value;\nconst openTagMatch = targetXml.match(/<w:p\\b[^>]*>/);\nlet pOpenTag = openTagMatch ? openTagMatch[0] : \"<w:p>\";\npOpenTag = pOpenTag\n  .replace(/\\s+w14:paraId=\"[^\"]*\"/, \"\")\n  .replace(/\\s+w14:textId=\"[^\"]*\"/, \"\");\n\n// Replace the whole paragraph content (the field run sequence + bookmark)\n// with literal text runs, keeping the bookmark in the same spot so the\n// \"_GoBack\" bookmark keeps pointing at the same logical location.\nconst contentRange = targetParagraph.getRange(\"Content\");\n\nconst newOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  pOpenTag +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>m</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  \"<w:r><w:t>'</w:t></w:r>\" +\n  '<w:r><w:t>doc.html</w:t></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  \"<w:r><w:t>'.fromHTMLURI()</w:t></w:r>\" +\n  '<w:r><w:t xml:space=\"preserve\">}</w:t></w:r>' +\n  '</w:p>' +\n  '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\ncontentRange.insertOoxml(newOoxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# The document contains a paragraph holding a Word FIELD (fldChar begin /\n# instrText tokens for \" m:'doc.html'.fromHTMLURI() \" / fldChar end).\n# This script rewrites that field into plain literal text runs that spell\n# out the same token stream wrapped in literal \"{\" and \"}\" characters\n# instead of the field delimiters (i.e. \"{m:'doc.html'.fromHTMLURI()}\"),\n# while leaving the bookmark (\"_GoBack\") that sits between the \"doc.html\"\n# and \"'.fromHTMLURI()\" tokens untouched.\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that contains the field (fldChar begin/end + instrText\n# runs). We detect it through its OOXML rather than assuming a fixed index,\n# so the script is resilient to minor document changes.\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $xml = $p.Range.WordOpenXML\n    if (($xml -like \"*instrText*\") -and ($xml -like \"*fldChar*\")) {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not locate the paragraph containing the field to rewrite.\"\n}\n\n$targetRange = $target.Range\n$targetXml = $targetRange.WordOpenXML\n\n# Recover the original <w:p ...> opening tag attributes (rsid bookkeeping,\n# etc.) so the rewritten paragraph keeps them instead of becoming a bare\n# <w:p>. The w14:paraId/w14:textId attributes are synthesized by the OOXML\n# exporter itself (not present in the source document), so they are\n# stripped back out.\n$pOpenTag = \"<w:p>\"\nif ($targetXml -match '<w:p\\b[^>]*>') {\n    $pOpenTag = $Matches[0]\n}\n$pOpenTag = $pOpenTag -replace ' w14:paraId=\"[^\"]*\"', ''\n$pOpenTag = $pOpenTag -replace ' w14:textId=\"[^\"]*\"', ''\n\n# Replace the whole paragraph content (the field run sequence + bookmark)\n# with literal text runs, keeping the bookmark in the same spot so the\n# \"_GoBack\" bookmark keeps pointing at the same logical location.\n$newOoxml = @\"\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          $pOpenTag\n            <w:r><w:t>{</w:t></w:r>\n            <w:r><w:t>m</w:t></w:r>\n            <w:r><w:t>:</w:t></w:r>\n            <w:r><w:t>'</w:t></w:r>\n            <w:r><w:t>doc.html</w:t></w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n            <w:r><w:t>'.fromHTMLURI()</w:t></w:r>\n            <w:r><w:t xml:space=\"preserve\">}</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n\"@\n\n$targetRange.InsertXML($newOoxml)\n"}
